$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# Price cells are forced to Text format so strings such as "37.172.18" or
# "0.0770" are preserved exactly instead of being parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.172.18"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.013.72"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.23"
$ws.Range("E5").Value = "  +4.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -1.24%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.36"
$ws.Range("E8").Value = "  -7.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  -2.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0770"
$ws.Range("E10").Value = "  -5.27%  "

$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.38"
$ws.Range("E12").Value = "  -4.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.306.43"
$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.39"
$ws.Range("E14").Value = "  -2.66%  "

$ws.Range("E15").Value = "  -6.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.22"
$ws.Range("E16").Value = "  -4.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.001.25"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.079.20"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.25"
$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0836"
$ws.Range("E20").Value = "  -3.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "234.23"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.10"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.09"
$ws.Range("E26").Value = "  +0.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.97"
$ws.Range("E27").Value = "  -4.64%  "

$ws.Range("E28").Value = "  -1.46%  "

$ws.Range("E29").Value = "  -8.24%  "

$ws.Range("E30").Value = "  -2.87%  "

$ws.Range("E31").Value = "  -1.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.60"
$ws.Range("E32").Value = "  -3.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0642"
$ws.Range("E33").Value = "  -4.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.42"
$ws.Range("E34").Value = "  -1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("E35").Value = "  -5.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.53"
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.54"
$ws.Range("E39").Value = "  +3.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").Value = "  +0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.19"
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.443.27"
$ws.Range("E42").Value = "  +4.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0926"
$ws.Range("E43").Value = "  -4.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0211"
$ws.Range("E44").Value = "  -2.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.47"
$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.60"
$ws.Range("E46").Value = "  -8.17%  "

$ws.Range("E47").Value = "  -3.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.94"
$ws.Range("E48").Value = "  +2.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.96"
$ws.Range("E49").Value = "  -6.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.198.12"
$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("E51").Value = "  -7.82%  "
